$wb = $excel.ActiveWorkbook

$oldGuid = "62d12a74-00a3-4570-8909-480b2f895fd0"
$newGuid = "1975f8f9-9818-43c9-9080-7f9bff8450b6"
$oldHash = "2bd50c29398a6731780edcf4858ade3d3a3f8dba"
$newHash = "3c513ebd15a8c702535255a83f402e99e053f38b"

function Update-Hyperlink($ws, $a1, $newText) {
    # Re-point the single hyperlink anchored at $a1 so its TextToDisplay
    # becomes $newText while leaving its Address/SubAddress/ScreenTip (and
    # every other hyperlink on the sheet) exactly as they were. Editing
    # TextToDisplay in place on an item pulled from Hyperlinks.Item(n)
    # stacks a duplicate link instead of updating it, so capture the
    # existing target first, delete just this one link, then re-add it.
    $target = $null
    foreach ($h in @($ws.Hyperlinks)) {
        if ($h.Range.Address() -eq $a1) {
            $target = $h
            break
        }
    }
    $addr = $target.Address
    $subAddr = $target.SubAddress
    $tip = $target.ScreenTip
    $target.Delete()
    $ws.Hyperlinks.Add($ws.Range($a1), $addr, $subAddr, $tip, $newText) | Out-Null
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("D2").Value = "2016-46-19 02:46:42"
Update-Hyperlink $wsOverview '$A$2' "$newGuid.md"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("D2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-19 02:46:39"
Update-Hyperlink $wsZh '$A$2' "$newGuid.md"
Update-Hyperlink $wsZh '$D$2' "$newGuid.$newHash.zh-cn.xlf"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("D2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-19 02:46:42"
Update-Hyperlink $wsDe '$A$2' "$newGuid.md"
Update-Hyperlink $wsDe '$D$2' "$newGuid.$newHash.de-de.xlf"
